# Apply cryptos list price/volume updates (GitHub Actions sync)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.725.90'
$ws.Range('E2').Value = '  -0.17%  '
$ws.Range('D3').Value = '1.637.82'
$ws.Range('E3').Value = '  -0.70%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '217.94'
$ws.Range('E5').Value = '  +0.60%  '
$ws.Range('E6').Value = '  -1.35%  '
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('E8').Value = '  -1.27%  '
$ws.Range('E9').Value = '  -1.29%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '18.92'
$ws.Range('E10').Value = '  -1.68%  '
$ws.Range('E11').Value = '  +0.01%  '
$ws.Range('D12').Value = '1.868.70'
$ws.Range('E12').Value = '  -0.56%  '
$ws.Range('D13').Value = '1.640.19'
$ws.Range('E13').Value = '  -0.60%  '
$ws.Range('E14').Value = '  -2.03%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.522'
$ws.Range('E15').Value = '  -2.10%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '64.03'
$ws.Range('E16').Value = '  -2.04%  '
$ws.Range('D17').Value = '26.723.56'
$ws.Range('E17').Value = '  -0.22%  '
$ws.Range('E18').Value = '  -3.30%  '
$ws.Range('E19').Value = '  +0.07%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '209.02'
$ws.Range('E20').Value = '  -3.94%  '
$ws.Range('E21').Value = '  -1.36%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.16'
$ws.Range('E22').Value = '  -1.76%  '
$ws.Range('E23').Value = '  -7.26%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '9.16'
$ws.Range('E24').Value = '  -3.46%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '147.00'
$ws.Range('E25').Value = '  -0.13%  '
$ws.Range('E26').Value = '  +0.43%  '
$ws.Range('E27').Value = '  -2.46%  '
$ws.Range('E28').Value = '  -2.28%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '15.48'
$ws.Range('E29').Value = '  -1.89%  '
$ws.Range('E30').Value = '  -3.93%  '
$ws.Range('E31').Value = '  +0.66%  '
$ws.Range('E32').Value = '  -1.19%  '
$ws.Range('E33').Value = '  -2.32%  '
$ws.Range('D34').Value = '1.263.00'
$ws.Range('E34').Value = '  -1.44%  '
$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.51'
$ws.Range('E35').Value = '  -2.21%  '
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.45'
$ws.Range('E36').Value = '  -0.03%  '
$ws.Range('E37').Value = '  -3.22%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.520'
$ws.Range('E38').Value = '  -3.34%  '
$ws.Range('E39').Value = '  -0.02%  '
$ws.Range('E40').Value = '  -3.75%  '
$ws.Range('E41').Value = '  -2.15%  '
$ws.Range('E42').Value = '  -2.98%  '
$ws.Range('D43').Value = '1.780.28'
$ws.Range('E43').Value = '  -0.55%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '5.24'
$ws.Range('E44').Value = '  -3.76%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '90.88'
$ws.Range('E45').Value = '  -1.26%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '59.70'
$ws.Range('E46').Value = '  -0.16%  '
$ws.Range('E47').Value = '  -2.57%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.0₆0102'
$ws.Range('E48').Value = '  -1.87%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0519'
$ws.Range('E49').Value = '  +0.74%  '
$ws.Range('B50').Value = 'USDD'
$ws.Range('C50').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.01'
$ws.Range('E50').Value = '  +0.11%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.407'
$ws.Range('E51').Value = '  -0.48%  '
